# Weekly update: add the newest reporting date (2023-11-28, Excel serial 45258)
# as a new Primera/Segunda pair at the top of the data block (rows 436-437),
# pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above row 436 (calling Insert twice on the same row
# index shifts the data down one row at a time).
$ws.Rows.Item(436).Insert()
$ws.Rows.Item(436).Insert()

# Row 436 - Calidad "Primera"
$ws.Cells.Item(436, 1).Value2  = 11
$ws.Cells.Item(436, 2).Value2  = 'Vega Monumental Concepción'
$ws.Cells.Item(436, 3).Value2  = 'Bíobío'
$ws.Cells.Item(436, 4).Value2  = 45258
$ws.Cells.Item(436, 5).Value2  = 8
$ws.Cells.Item(436, 6).Value2  = 100112008
$ws.Cells.Item(436, 7).Value2  = 'Coliflor'
$ws.Cells.Item(436, 8).Value2  = 'Sin especificar'
$ws.Cells.Item(436, 9).Value2  = 'Primera'
$ws.Cells.Item(436, 10).Value2 = 1000
$ws.Cells.Item(436, 11).Value2 = 1000
$ws.Cells.Item(436, 12).Value2 = 1000
$ws.Cells.Item(436, 13).Value2 = 1000
$ws.Cells.Item(436, 14).Value2 = '$/unidad'
$ws.Cells.Item(436, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(436, 16).Value2 = 1000
$ws.Cells.Item(436, 17).Value2 = 1
$ws.Cells.Item(436, 18).Value2 = 'Hortaliza'

# Row 437 - Calidad "Segunda"
$ws.Cells.Item(437, 1).Value2  = 11
$ws.Cells.Item(437, 2).Value2  = 'Vega Monumental Concepción'
$ws.Cells.Item(437, 3).Value2  = 'Bíobío'
$ws.Cells.Item(437, 4).Value2  = 45258
$ws.Cells.Item(437, 5).Value2  = 8
$ws.Cells.Item(437, 6).Value2  = 100112008
$ws.Cells.Item(437, 7).Value2  = 'Coliflor'
$ws.Cells.Item(437, 8).Value2  = 'Sin especificar'
$ws.Cells.Item(437, 9).Value2  = 'Segunda'
$ws.Cells.Item(437, 10).Value2 = 1000
$ws.Cells.Item(437, 11).Value2 = 1000
$ws.Cells.Item(437, 12).Value2 = 1000
$ws.Cells.Item(437, 13).Value2 = 1000
$ws.Cells.Item(437, 14).Value2 = '$/unidad'
$ws.Cells.Item(437, 15).Value2 = 'Región Metropolitana'
$ws.Cells.Item(437, 16).Value2 = 1000
$ws.Cells.Item(437, 17).Value2 = 1
$ws.Cells.Item(437, 18).Value2 = 'Hortaliza'

# Match the date-column number format used by the rest of column D.
$ws.Range("D436:D437").NumberFormat = "YYYY-MM-DD HH:MM:SS"
